# Add 4 new data rows (9-12) to "Лист1", mirroring the existing table
# (rows 2-8: A=длина ряда, B=размер элемента, C=длина подпоследовательности,
#  D..G computed via formulas), as part of "optimizations for paral & vect".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number format strings copied verbatim (with literal-space backslash escapes)
# from the existing financial formats used by columns A/D/E/F and G, so that
# the new cells reuse the very same numFmt definitions already in the file.
$rub = [char]0x20BD
$fmt0 = "_-* #,##0\ _${rub}_-;\-* #,##0\ _${rub}_-;_-* ""-""??\ _${rub}_-;_-@_-"
$fmt2 = "_-* #,##0.00\ _${rub}_-;\-* #,##0.00\ _${rub}_-;_-* ""-""??\ _${rub}_-;_-@_-"

$newRows = @(
    @{ Row = 9;  A = 90000; B = 4; C = 1024 },
    @{ Row = 10; A = 90000; B = 4; C = 128 },
    @{ Row = 11; A = 65000; B = 4; C = 1024 },
    @{ Row = 12; A = 65000; B = 4; C = 128 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)
    $cellE = $ws.Cells.Item($r, 5)
    $cellF = $ws.Cells.Item($r, 6)
    $cellG = $ws.Cells.Item($r, 7)

    $cellA.Value = $item.A
    $cellB.Value = $item.B
    $cellC.Value = $item.C

    $cellD.Formula = "=A$r-C$r+1"
    $cellE.Formula = "=D$r*D$r"
    $cellF.Formula = "=C$r*D$r"
    $cellG.Formula = "=(E$r+F$r)*B$r/1000000000"

    # Formatting to match rows 2-8: column A & G use the financial formats,
    # B/C/D/E/F share the plain "vertical top" alignment used throughout.
    $cellA.NumberFormat = $fmt0
    $cellG.NumberFormat = $fmt2

    $cellA.VerticalAlignment = -4160
    $cellB.VerticalAlignment = -4160
    $cellC.VerticalAlignment = -4160
    $cellD.VerticalAlignment = -4160
    $cellE.VerticalAlignment = -4160
    $cellF.VerticalAlignment = -4160
    $cellG.VerticalAlignment = -4160

    $cellD.NumberFormat = $fmt0
    $cellE.NumberFormat = $fmt0
    $cellF.NumberFormat = $fmt0
}

# Match the selection left behind by the author (selecting the next empty row).
$ws.Range("A13").Select() | Out-Null
